# Insert a new top entry ("Puerto Rico", id 70) above the existing ranking
# table, pushing every existing row (and its hyperlink) down by one row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new blank row at row 2; everything below (values, styles,
#    hyperlinks anchor cells) shifts down by one row.
$ws.Rows("2:2").Insert()

# 2. Fill in the new row with the Puerto Rico entry (rank/id 70).
$ws.Range("A2").Value2 = 70
$ws.Range("B2").Value2 = "pr"
$ws.Range("C2").Value2 = "Puerto Rico"

# 3. The underlying engine does not shift the hyperlink anchors together
#    with the row insert, so rebuild the hyperlink collection from scratch
#    at the correct (now shifted-down-by-one) rows: old A7..A70 -> A8..A71.
$ws.Cells.Hyperlinks.Delete()
for ($display = 64; $display -ge 1; $display--) {
    $row = 72 - $display
    $cell = $ws.Cells.Item($row, 1)
    $url = "http://127.0.0.1:8050/en/admin/avi/countrydim/" + $display + "/change/"
    $ws.Hyperlinks.Add($cell, $url, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, [string]$display)
}

# 4. Match the recorded selection state after the edit.
$ws.Range("A3").Select()
